$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(33, 8).Value = 1309.25
$ws.Cells.Item(33, 9).Value = 1604.5
$ws.Cells.Item(33, 10).Value = 423.5
$ws.Cells.Item(33, 11).Value = 1604.5
$ws.Cells.Item(33, 12).Value = 423.5
$ws.Cells.Item(33, 13).Value = -1375.5
$ws.Cells.Item(33, 14).Value = -881.5
$ws.Cells.Item(69, 8).Value = 20389.389
$ws.Cells.Item(69, 10).Value = 20588.941
$ws.Cells.Item(69, 12).Value = 61766.823
$ws.Cells.Item(69, 14).Value = -63514.823
$ws.Cells.Item(72, 8).Value = 20389.389
$ws.Cells.Item(72, 10).Value = 20588.941
$ws.Cells.Item(72, 12).Value = 185300.469
$ws.Cells.Item(72, 14).Value = -194036.469
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(99, 8).Value = 378
$ws.Cells.Item(99, 9).Value = 282.8
$ws.Cells.Item(99, 10).Value = 536.6667
$ws.Cells.Item(99, 11).Value = 848.4000000000001
$ws.Cells.Item(99, 12).Value = 1610.0001
$ws.Cells.Item(99, 13).Value = 649.5999999999999
$ws.Cells.Item(99, 14).Value = -4606.0001
$ws.Cells.Item(101, 8).Value = 16668115
$ws.Cells.Item(101, 9).Value = 20000738
$ws.Cells.Item(101, 11).Value = 60002214
$ws.Cells.Item(101, 13).Value = -60000592
$ws.Cells.Item(134, 8).Value = 45800
$ws.Cells.Item(134, 10).Value = 45800
$ws.Cells.Item(134, 12).Value = 45800
$ws.Cells.Item(134, 14).Value = -55940
$ws.Cells.Item(137, 8).Value = 1525.35
$ws.Cells.Item(137, 9).Value = 750.5714
$ws.Cells.Item(137, 11).Value = 2251.7142
$ws.Cells.Item(137, 13).Value = 298.2857999999997

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 1880206.6
$ws.Cells.Item(32, 9).Value = 2123394
$ws.Cells.Item(32, 11).Value = 2123394
$ws.Cells.Item(32, 13).Value = -2123107
$ws.Cells.Item(43, 8).Value = 44950
$ws.Cells.Item(43, 9).Value = 44950
$ws.Cells.Item(43, 11).Value = 44950
$ws.Cells.Item(43, 13).Value = -44637
$ws.Cells.Item(88, 8).Value = 3251.25
$ws.Cells.Item(88, 9).Value = 2918.3333
$ws.Cells.Item(88, 10).Value = 4250
$ws.Cells.Item(88, 11).Value = 2918.3333
$ws.Cells.Item(88, 12).Value = 4250
$ws.Cells.Item(88, 13).Value = -2512.3333
$ws.Cells.Item(88, 14).Value = -5062
$ws.Cells.Item(91, 8).Value = 3251.25
$ws.Cells.Item(91, 9).Value = 2918.3333
$ws.Cells.Item(91, 10).Value = 4250
$ws.Cells.Item(91, 11).Value = 2918.3333
$ws.Cells.Item(91, 12).Value = 4250
$ws.Cells.Item(91, 13).Value = -1514.3333
$ws.Cells.Item(91, 14).Value = -7058

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 1001.8
$ws.Cells.Item(20, 9).Value = 1016.6667
$ws.Cells.Item(20, 11).Value = 1016.6667
$ws.Cells.Item(20, 13).Value = -769.6667
$ws.Cells.Item(86, 8).Value = 2060
$ws.Cells.Item(86, 9).Value = 1913.3334
$ws.Cells.Item(86, 11).Value = 1913.3334
$ws.Cells.Item(86, 13).Value = -790.3334
$ws.Cells.Item(89, 8).Value = 2060
$ws.Cells.Item(89, 9).Value = 1913.3334
$ws.Cells.Item(89, 11).Value = 9566.666999999999
$ws.Cells.Item(89, 13).Value = -3950.666999999999

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 8).Value = 1000
$ws.Cells.Item(3, 10).Value = 1000
$ws.Cells.Item(3, 12).Value = 1000
$ws.Cells.Item(3, 14).Value = -1226
$ws.Cells.Item(16, 8).Value = 2432.3572
$ws.Cells.Item(16, 9).Value = 2636.1
$ws.Cells.Item(16, 11).Value = 2636.1
$ws.Cells.Item(16, 13).Value = -2349.1
$ws.Cells.Item(22, 8).Value = 594.3077
$ws.Cells.Item(22, 9).Value = 546.375
$ws.Cells.Item(22, 10).Value = 671
$ws.Cells.Item(22, 11).Value = 546.375
$ws.Cells.Item(22, 12).Value = 671
$ws.Cells.Item(22, 13).Value = -196.375
$ws.Cells.Item(22, 14).Value = -1371
$ws.Cells.Item(62, 8).Value = 3397.5
$ws.Cells.Item(62, 9).Value = 3800
$ws.Cells.Item(62, 11).Value = 3800
$ws.Cells.Item(62, 13).Value = -3176
$ws.Cells.Item(65, 8).Value = 3397.5
$ws.Cells.Item(65, 9).Value = 3800
$ws.Cells.Item(65, 11).Value = 19000
$ws.Cells.Item(65, 13).Value = -15880
$ws.Cells.Item(105, 8).Value = 3164.9333
$ws.Cells.Item(105, 9).Value = 2369.2856
$ws.Cells.Item(105, 10).Value = 3861.125
$ws.Cells.Item(105, 11).Value = 2369.2856
$ws.Cells.Item(105, 12).Value = 3861.125
$ws.Cells.Item(105, 13).Value = -622.2856000000002
$ws.Cells.Item(105, 14).Value = -7355.125
$ws.Cells.Item(113, 8).Value = 2432.3572
$ws.Cells.Item(113, 9).Value = 2636.1
$ws.Cells.Item(113, 11).Value = 2636.1
$ws.Cells.Item(113, 13).Value = -466.0999999999999
$ws.Cells.Item(134, 8).Value = 1891.8572
$ws.Cells.Item(134, 9).Value = 1891.8572
$ws.Cells.Item(134, 11).Value = 5675.571599999999
$ws.Cells.Item(134, 13).Value = -3140.571599999999
$ws.Cells.Item(141, 8).Value = 34989.547
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 34989.547
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 34989.547
$ws.Cells.Item(141, 13).ClearContents()
$ws.Cells.Item(141, 14).Value = -45349.547

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(70, 8).Value = 9772.700000000001
$ws.Cells.Item(70, 9).Value = 1545.8
$ws.Cells.Item(70, 11).Value = 4637.4
$ws.Cells.Item(70, 13).Value = -4322.4
$ws.Cells.Item(73, 8).Value = 9772.700000000001
$ws.Cells.Item(73, 9).Value = 1545.8
$ws.Cells.Item(73, 11).Value = 4637.4
$ws.Cells.Item(73, 13).Value = -3545.4
$ws.Cells.Item(80, 8).Value = 7000
$ws.Cells.Item(80, 9).Value = 10000
$ws.Cells.Item(80, 10).Value = 4000
$ws.Cells.Item(80, 11).Value = 30000
$ws.Cells.Item(80, 12).Value = 12000
$ws.Cells.Item(80, 13).Value = -29064
$ws.Cells.Item(80, 14).Value = -13872
$ws.Cells.Item(83, 8).Value = 7000
$ws.Cells.Item(83, 9).Value = 10000
$ws.Cells.Item(83, 10).Value = 4000
$ws.Cells.Item(83, 11).Value = 90000
$ws.Cells.Item(83, 12).Value = 36000
$ws.Cells.Item(83, 13).Value = -85320
$ws.Cells.Item(83, 14).Value = -45360
$ws.Cells.Item(92, 8).Value = 670.85
$ws.Cells.Item(92, 9).Value = 566
$ws.Cells.Item(92, 10).Value = 689.35297
$ws.Cells.Item(92, 11).Value = 1698
$ws.Cells.Item(92, 12).Value = 2068.05891
$ws.Cells.Item(92, 13).Value = -450
$ws.Cells.Item(92, 14).Value = -4564.05891
$ws.Cells.Item(141, 8).Value = 7451
$ws.Cells.Item(141, 9).Value = 7757.375
$ws.Cells.Item(141, 11).Value = 23272.125
$ws.Cells.Item(141, 13).Value = -18092.125

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(80, 8).Value = 3987.889
$ws.Cells.Item(80, 9).Value = 3932.3333
$ws.Cells.Item(80, 10).Value = 4015.6667
$ws.Cells.Item(80, 11).Value = 3932.3333
$ws.Cells.Item(80, 12).Value = 4015.6667
$ws.Cells.Item(80, 13).Value = -2934.3333
$ws.Cells.Item(80, 14).Value = -6011.6667
$ws.Cells.Item(83, 8).Value = 3987.889
$ws.Cells.Item(83, 9).Value = 3932.3333
$ws.Cells.Item(83, 10).Value = 4015.6667
$ws.Cells.Item(83, 11).Value = 19661.6665
$ws.Cells.Item(83, 12).Value = 20078.3335
$ws.Cells.Item(83, 13).Value = -14669.6665
$ws.Cells.Item(83, 14).Value = -30062.3335
$ws.Cells.Item(122, 8).Value = 2297.6667
$ws.Cells.Item(122, 9).Value = 1957.2
$ws.Cells.Item(122, 11).Value = 5871.6
$ws.Cells.Item(122, 13).Value = -3421.6

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 833.1111
$ws.Cells.Item(22, 9).Value = 883
$ws.Cells.Item(22, 10).Value = 733.3333
$ws.Cells.Item(22, 11).Value = 883
$ws.Cells.Item(22, 12).Value = 733.3333
$ws.Cells.Item(22, 13).Value = -588
$ws.Cells.Item(22, 14).Value = -1323.3333
$ws.Cells.Item(27, 8).Value = 833.1111
$ws.Cells.Item(27, 9).Value = 883
$ws.Cells.Item(27, 10).Value = 733.3333
$ws.Cells.Item(27, 11).Value = 883
$ws.Cells.Item(27, 12).Value = 733.3333
$ws.Cells.Item(27, 13).Value = -776
$ws.Cells.Item(27, 14).Value = -947.3333
$ws.Cells.Item(39, 8).Value = 19999
$ws.Cells.Item(39, 10).Value = 19999
$ws.Cells.Item(39, 12).Value = 19999
$ws.Cells.Item(39, 14).Value = -20919
$ws.Cells.Item(132, 8).Value = 6751.4
$ws.Cells.Item(132, 9).Value = 7568.727
$ws.Cells.Item(132, 11).Value = 22706.181
$ws.Cells.Item(132, 13).Value = -20176.181

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(54, 8).Value = 22333
$ws.Cells.Item(54, 10).Value = 44999
$ws.Cells.Item(54, 12).Value = 44999
$ws.Cells.Item(54, 14).Value = -46039
$ws.Cells.Item(70, 8).Value = 90095
$ws.Cells.Item(70, 9).Value = 90095
$ws.Cells.Item(70, 11).Value = 90095
$ws.Cells.Item(70, 13).Value = -89780
$ws.Cells.Item(73, 8).Value = 90095
$ws.Cells.Item(73, 9).Value = 90095
$ws.Cells.Item(73, 11).Value = 90095
$ws.Cells.Item(73, 13).Value = -89003
$ws.Cells.Item(136, 8).Value = 2250.9565
$ws.Cells.Item(136, 9).Value = 2110.1667
$ws.Cells.Item(136, 10).Value = 2757.8
$ws.Cells.Item(136, 11).Value = 6330.500100000001
$ws.Cells.Item(136, 12).Value = 8273.400000000001
$ws.Cells.Item(136, 13).Value = -3780.500100000001
$ws.Cells.Item(136, 14).Value = -13373.4
